$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 86 — "lang_select_image" key/vi/en triple.
# Copy formatting from row 83 (A = key style, B/C = value style) so the
# new row reuses the existing style entries instead of minting new ones.
$ws.Range("A83:C83").Copy()
$ws.Range("A86:C86").PasteSpecial(-4122)
$ws.Range("A86").Value = "lang_select_image"
$ws.Range("B86").Value = "Chọn Ảnh"
$ws.Range("C86").Value = "Select Image"

# New row 87 — "lang_full_name" key/vi/en triple.
# Copy formatting from row 84, whose three cells all share the plain style.
$ws.Range("A84:C84").Copy()
$ws.Range("A87:C87").PasteSpecial(-4122)
$ws.Range("A87").Value = "lang_full_name"
$ws.Range("B87").Value = "Họ & Tên"
$ws.Range("C87").Value = "Full Name"

# Match the author's final selection in the saved workbook.
$ws.Range("C87").Select()
